# fix(builder): ensure output folder exists before saving .pptx presentation
#
# Rewrites the bullet-point slide contents (titles, bullet lines, and the
# bold "takeaway" line) across the deck. Uses character-range replacement
# (rather than whole-TextRange reassignment) so existing per-paragraph
# formatting (bullet sz=1800, bold pink takeaway line) is preserved and no
# stray runs are introduced.

function Set-ParaText($tr, $idx, $newText) {
    # Replace the text of paragraph number $idx (1-based) in TextRange $tr,
    # keeping that paragraph's own formatting (pPr/defRPr) untouched and
    # leaving exactly one run behind.
    $para = $tr.Paragraphs($idx, 1)
    $sub = $tr.Characters($para.Start, $para.Length - 1)
    $sub.Text = $newText
}

function Set-ShapeText($sh, $newText) {
    # Replace the (single-paragraph) text of a shape's TextRange in place.
    $tr = $sh.TextFrame.TextRange
    $len = $tr.Text.Length
    $sub = $tr.Characters(1, $len)
    $sub.Text = $newText
}

$p = $ppt.ActivePresentation

# --- Slide 2: "Welcome to Our Feline Friends" -> "Welcome to the World of Cats!"
$s = $p.Slides.Item(2)
Set-ShapeText $s.Shapes.Item(1) "Welcome to the World of Cats!"
$tr = $s.Shapes.Item(2).TextFrame.TextRange
$n = $tr.Paragraphs().Count
Set-ParaText $tr $n "Did you know? A group of cats is called a clowder!"

# --- Slide 3: "Why Cats Make Great House Pets" -> "Why Cats Make Great Pets"
$s = $p.Slides.Item(3)
Set-ShapeText $s.Shapes.Item(1) "Why Cats Make Great Pets"
$tr = $s.Shapes.Item(2).TextFrame.TextRange
Set-ParaText $tr 1 "Independent and low maintenance"
Set-ParaText $tr 2 "Provide companionship"
Set-ParaText $tr 3 "Natural hunters"
$n = $tr.Paragraphs().Count
Set-ParaText $tr $n "Cats are perfect for those who appreciate independence and companionship."

# --- Slide 5: "The Science Behind the Purr" -> "Case Study: Cats in Therapy"
$s = $p.Slides.Item(5)
Set-ShapeText $s.Shapes.Item(1) "Case Study: Cats in Therapy"
$tr = $s.Shapes.Item(2).TextFrame.TextRange
Set-ParaText $tr 1 "Therapy cats in hospitals"
Set-ParaText $tr 2 "Emotional support animals"
Set-ParaText $tr 3 "Positive impact on mental health"
$n = $tr.Paragraphs().Count
Set-ParaText $tr $n "Cats are increasingly used in therapeutic settings to aid mental health."

# --- Slide 6: "Case Study: Cats in Therapy" -> "Unique Cat Behaviors"
$s = $p.Slides.Item(6)
Set-ShapeText $s.Shapes.Item(1) "Unique Cat Behaviors"
$tr = $s.Shapes.Item(2).TextFrame.TextRange
Set-ParaText $tr 1 "Kneading with paws"
Set-ParaText $tr 2 "Zoomies: sudden energy bursts"
Set-ParaText $tr 3 "Hiding in small spaces"
$n = $tr.Paragraphs().Count
Set-ParaText $tr $n "These behaviors can be entertaining and comforting."

# --- Slide 7: "Cats' Unique Behaviors" -> "Cats' Communication Methods"
$s = $p.Slides.Item(7)
Set-ShapeText $s.Shapes.Item(1) "Cats' Communication Methods"
$tr = $s.Shapes.Item(2).TextFrame.TextRange
Set-ParaText $tr 1 "Vocalizations: meowing and purring"
Set-ParaText $tr 2 "Body language"
Set-ParaText $tr 3 "Tail movements"
$n = $tr.Paragraphs().Count
Set-ParaText $tr $n "Understanding cat communication can enhance the human-cat bond."

# --- Slide 8: "Step-by-Step: Introducing a Cat to Your Home" -> "Historical Significance of Cats"
#     (4 bullets -> 3 bullets: drop the 4th one)
$s = $p.Slides.Item(8)
Set-ShapeText $s.Shapes.Item(1) "Historical Significance of Cats"
$tr = $s.Shapes.Item(2).TextFrame.TextRange
Set-ParaText $tr 1 "Domesticated around 7500 BC"
Set-ParaText $tr 2 "Worshipped in ancient Egypt"
Set-ParaText $tr 3 "Mummification practices"
$tr.Paragraphs(4, 1).Delete()
$n = $tr.Paragraphs().Count
Set-ParaText $tr $n "Cats have been valued companions throughout history."

# --- Slide 9: "Historical Significance of Cats" -> "Cats' Physical Abilities"
$s = $p.Slides.Item(9)
Set-ShapeText $s.Shapes.Item(1) "Cats' Physical Abilities"
$tr = $s.Shapes.Item(2).TextFrame.TextRange
Set-ParaText $tr 1 "Excellent night vision"
Set-ParaText $tr 2 "Hear ultrasonic sounds"
Set-ParaText $tr 3 "Agility and hunting skills"
$n = $tr.Paragraphs().Count
Set-ParaText $tr $n "Cats are equipped with remarkable physical traits."

# --- Slide 10: "Fun Facts About Cats" -> "Step-by-Step: Introducing a Cat to Your Home"
#     (3 bullets -> 4 bullets: add a 4th one)
$s = $p.Slides.Item(10)
Set-ShapeText $s.Shapes.Item(1) "Step-by-Step: Introducing a Cat to Your Home"
$tr = $s.Shapes.Item(2).TextFrame.TextRange
Set-ParaText $tr 1 ">> Prepare a safe space"
Set-ParaText $tr 2 ">> Gradually introduce to family members"
Set-ParaText $tr 3 ">> Establish a feeding routine"
[void]$tr.Paragraphs(3, 1).InsertAfter("`r>> Provide toys and scratching posts")
$n = $tr.Paragraphs().Count
Set-ParaText $tr $n "A smooth introduction ensures a happy home for your cat."

# --- Slide 11: "Key Takeaways" (title unchanged)
$s = $p.Slides.Item(11)
$tr = $s.Shapes.Item(2).TextFrame.TextRange
Set-ParaText $tr 1 "Cats offer companionship and anxiety relief"
Set-ParaText $tr 2 "Unique behaviors and communication"
Set-ParaText $tr 3 "Historical and cultural significance"
$n = $tr.Paragraphs().Count
Set-ParaText $tr $n "Cats are more than just pets; they are companions with rich histories."

# --- Slide 12: "Conclusion" (title unchanged)
$s = $p.Slides.Item(12)
$tr = $s.Shapes.Item(2).TextFrame.TextRange
Set-ParaText $tr 1 "Cats are ideal for various lifestyles"
Set-ParaText $tr 2 "They provide emotional support and entertainment"
$n = $tr.Paragraphs().Count
Set-ParaText $tr $n "Consider adopting a cat for a fulfilling companionship."

# --- Slide 13: "Thank You!" (title & bottom line unchanged, only 2nd bullet changes)
$s = $p.Slides.Item(13)
$tr = $s.Shapes.Item(2).TextFrame.TextRange
Set-ParaText $tr 2 "Let's discuss the joys of having cats!"
